$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Fuel"
$ws.Range("B2").Value = "Oil"

$ws.Range("A3").Select()
